# Q3 Update - 2025
# Applies the shared-string content updates described by the diff:
#   1. short-url column (B2:B111) "cALvB0" -> "o7PpJ1" for every data row
#   2. refugees value for row 110 (N110) "8955" -> "8505"
#   3. stateless value for row 111 (S111) "22496" -> "20000"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) short-url column: replace the shared text across all data rows
$null = $ws.Range("B2:B111").Replace("cALvB0", "o7PpJ1")

# 2) refugees (row 110) changes from 8955 to 8505 - keep it stored as text,
#    matching the rest of the text-typed numeric data in this sheet.
$ws.Range("N110").NumberFormat = "@"
$ws.Range("N110").Value = "8505"

# 3) stateless (row 111) changes from 22496 to 20000 - also stored as text.
$ws.Range("S111").NumberFormat = "@"
$ws.Range("S111").Value = "20000"
